$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.692.14"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "2.919.10"
$ws.Range("E3").Value = "  +1.30%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "354.61"
$ws.Range("E5").Value = "  +0.69%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "110.41"
$ws.Range("E6").Value = "  -1.01%  "

$ws.Range("E7").Value = "  +1.51%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.628"
$ws.Range("E9").Value = "  +1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.25"
$ws.Range("E10").Value = "  -2.07%  "

$ws.Range("E11").Value = "  +3.03%  "

$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("E13").Value = "  -1.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.90"
$ws.Range("E14").Value = "  +1.10%  "

$ws.Range("D15").Value = "3.378.65"
$ws.Range("E15").Value = "  +1.29%  "

$ws.Range("D16").Value = "2.921.08"
$ws.Range("E16").Value = "  +1.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.977"
$ws.Range("E17").Value = "  -1.78%  "

$ws.Range("D18").Value = "51.739.19"
$ws.Range("E18").Value = "  -0.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  -2.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.26"
$ws.Range("E20").Value = "  -2.74%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.96"
$ws.Range("E21").Value = "  +0.50%  "

$ws.Range("D22").Value = "0.0₃0980"
$ws.Range("E22").Value = "  -0.17%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.81"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.99"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.82"
$ws.Range("E25").Value = "  +1.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.06"
$ws.Range("E27").Value = "  +2.82%  "

$ws.Range("E28").Value = "  +0.16%  "

$ws.Range("E29").Value = "  +15.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.107"
$ws.Range("E30").Value = "  +14.21%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "38.32"
$ws.Range("E32").Value = "  -1.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.07"
$ws.Range("E33").Value = "  +2.68%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "52.31"
$ws.Range("E34").Value = "  -1.91%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0441"
$ws.Range("E35").Value = "  -3.77%  "

$ws.Range("E36").Value = "  +0.06%  "

$ws.Range("E37").Value = "  -15.90%  "

$ws.Range("E38").Value = "  -2.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.43"
$ws.Range("E39").Value = "  -0.84%  "

$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("E41").Value = "  +3.80%  "

$ws.Range("E42").Value = "  +1.40%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "22.87"
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.84"
$ws.Range("E44").Value = "  -1.52%  "

$ws.Range("E45").Value = "  -2.52%  "

$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.43"
$ws.Range("E47").Value = "  -4.04%  "

$ws.Range("D48").Value = "2.131.79"
$ws.Range("E48").Value = "  -3.11%  "

$ws.Range("E49").Value = "  -7.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0332"
$ws.Range("E50").Value = "  +4.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "9.08"
$ws.Range("E51").Value = "  -0.46%  "
